$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (e.g. "359.31", "1.00"). Force them to keep a Text format while writing
# the value, then restore the default "Normal" style so the cell keeps no
# explicit style, exactly like the rest of the sheet.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D5').Value = '359.31'
$ws.Range('D6').Value = '109.81'
$ws.Range('D7').Value = '0.559'
$ws.Range('D9').Value = '0.590'
$ws.Range('D10').Value = '39.76'
$ws.Range('D11').Value = '0.0848'
$ws.Range('D13').Value = '19.37'
$ws.Range('D14').Value = '7.53'
$ws.Range('D17').Value = '0.941'
$ws.Range('D19').Value = '7.49'
$ws.Range('D20').Value = '3.09'
$ws.Range('D21').Value = '13.07'
$ws.Range('D23').Value = '70.06'
$ws.Range('D24').Value = '269.79'
$ws.Range('D26').Value = '26.48'
$ws.Range('D27').Value = '1.00'
$ws.Range('D28').Value = '0.164'
$ws.Range('D29').Value = '10.20'
$ws.Range('D30').Value = '2.15'
$ws.Range('D31').Value = '0.0469'
$ws.Range('D32').Value = '51.91'
$ws.Range('D33').Value = '33.77'
$ws.Range('D34').Value = '5.75'
$ws.Range('D35').Value = '0.0841'
$ws.Range('D36').Value = '5.21'
$ws.Range('D38').Value = '18.58'
$ws.Range('D39').Value = '3.20'
$ws.Range('D40').Value = '1.99'
$ws.Range('D44').Value = '119.43'
$ws.Range('D45').Value = '21.89'
$ws.Range('D47').Value = '3.23'
$ws.Range('D49').Value = '5.72'
$ws.Range('D50').Value = '0.947'
$ws.Range('D51').Value = '8.85'

$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'

# Remaining cells: plain text / already non-numeric-looking values
$ws.Range('D2').Value = '51.968.00'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.786.96'
$ws.Range('E3').Value = '  -1.96%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  -3.22%  '
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('E10').Value = '  -4.66%  '
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('E14').Value = '  -3.22%  '
$ws.Range('D15').Value = '3.220.36'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').Value = '2.788.85'
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('E17').Value = '  +4.28%  '
$ws.Range('D18').Value = '51.910.75'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  -2.43%  '
$ws.Range('E21').Value = '  -3.64%  '
$ws.Range('D22').Value = '0.0₃0976'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E24').Value = '  +0.99%  '
$ws.Range('E25').Value = '  -3.09%  '
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E28').Value = '  +16.90%  '
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('E30').Value = '  -2.82%  '
$ws.Range('B31').Value = 'VeChain'
$ws.Range('C31').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E31').Value = '  +5.41%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('E33').Value = '  -0.38%  '
$ws.Range('E34').Value = '  -2.57%  '
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E39').Value = '  -3.68%  '
$ws.Range('E40').Value = '  -4.38%  '
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('E42').Value = '  -1.56%  '
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('E44').Value = '  -6.12%  '
$ws.Range('E45').Value = '  -8.93%  '
$ws.Range('D46').Value = '2.079.75'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('E47').Value = '  -4.11%  '
$ws.Range('E48').Value = '  -1.82%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('E50').Value = '  -6.31%  '
$ws.Range('E51').Value = '  -1.78%  '
